$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I4").Value = 5.1
$ws.Range("F6").Value = 2.84
$ws.Range("G6").Value = 4.2
$ws.Range("H6").Value = 2.22
$ws.Range("J6").Value = 2.42
$ws.Range("K6").Value = 5.1
$ws.Range("P6").Value = 1.64
$ws.Range("F7").Value = 2.38
$ws.Range("I7").Value = 3.8
$ws.Range("J7").Value = 2.88
$ws.Range("F8").Value = 1.25
$ws.Range("I8").Value = 1.63
$ws.Range("P8").Value = 2.12
$ws.Range("Q8").Value = 1.7
$ws.Range("F9").Value = 1.93
$ws.Range("I9").Value = 4.1
$ws.Range("K9").Value = 4.6
$ws.Range("F10").Value = 5.2
$ws.Range("G10").Value = 5.4
$ws.Range("H10").Value = 1.7
$ws.Range("I10").Value = 1.71
$ws.Range("L10").Value = 1.27
$ws.Range("P10").Value = 2.32
$ws.Range("Q10").Value = 1.7
$ws.Range("R10").Value = 1.53
$ws.Range("S10").Value = 2.8
$ws.Range("T10").Value = 1.72
$ws.Range("U10").Value = 2.2
$ws.Range("V10").Value = 2.4
$ws.Range("AE10").Value = 18
$ws.Range("AO10").Value = 8.800000000000001
$ws.Range("L11").Value = 1.41
$ws.Range("U11").Value = 2.14
$ws.Range("AK11").Value = 27
$ws.Range("Q12").Value = 1.47
$ws.Range("AB12").Value = 17.5
$ws.Range("AC12").Value = 12
$ws.Range("AK12").Value = 20
$ws.Range("F13").Value = 1.75
$ws.Range("G13").Value = 1.77
$ws.Range("P13").Value = 2.2
$ws.Range("Q13").Value = 1.76
$ws.Range("T13").Value = 1.75
$ws.Range("V13").Value = 1.25
$ws.Range("W13").Value = 2.28
$ws.Range("AJ13").Value = 17
$ws.Range("AK13").Value = 16
$ws.Range("AN13").Value = 9
$ws.Range("AO13").Value = 75
$ws.Range("F14").Value = 1.9
$ws.Range("H14").Value = 4.1
$ws.Range("N14").Value = 1.76
$ws.Range("P14").Value = 1.76
